$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

$data = @(
    @("Column1", "Column2"),
    @("UiPath: Robotic Process Automation (RPA)", "COURSE"),
    @("Introducing Robotic Process Automation", "COURSE"),
    @("RPA: Automation Anywhere", "COURSE"),
    @("RPA, AI, and Cognitive Tech for Leaders", "COURSE"),
    @("UiPath Essential Training", "COURSE"),
    @("RPA: Automation Anywhere IQ Bot", "COURSE"),
    @("Robotic Process Automation: Tech Primer", "COURSE"),
    @("Blue Prism: Excel Automation", "COURSE"),
    @("Introducing Blue Prism", "COURSE"),
    @("Digital Technologies Case Studies: AI, IOT, Robotics, Blockchain", "COURSE")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
